$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set in_service (column E) to TRUE for rows 10-14 (extr1..extr5)
$ws.Range("E10:E14").Value = $true
